$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update summary figures (VALOR MORA total, Cant. Trabajadores)
# ---------------------------------------------------------------------------
$ws.Cells.Item(11, 5).Value = 583584    # E11 VALOR MORA
$ws.Cells.Item(13, 3).Value = 11        # C13 Cant. Trabajadores

# ---------------------------------------------------------------------------
# 2. Before removing the extra rows, copy the special "bottom of table"
#    formatting (thicker bottom border) from the current last data row (38)
#    onto what will become the new last data row (29) once rows 30:38 are
#    removed.
# ---------------------------------------------------------------------------
$ws.Range("B38:J38").Copy()
$ws.Range("B29:J29").PasteSpecial(-4122)   # xlPasteFormats

# ---------------------------------------------------------------------------
# 3. Write the new worker/account data into rows 16-29. Columns C (document
#    number) and E (period) hold numeric-looking values but are stored as
#    text in this template, so the number format is (re)applied as "@"
#    right before assigning the value to keep them as text without leaving
#    a stray quote-prefix flag on the cell.
# ---------------------------------------------------------------------------
$data = @(
    @("CC", "22494380",   "KELLY DEL PILAR CORDOBA MARTINEZ", "1701", 27578,  689455),
    @("CC", "9020222",    "JAISON BENAVIDES RODELO",          "1701", 160000, 737717),
    @("CC", "9022904",    "WILMER ANTONIO CRUZ MORENO",       "1701", 27578,  737717),
    @("CC", "3871999",    "ADOLFO RAMON VARGAS NAVARRO",      "1701", 27578,  737717),
    @("CC", "1052971182", "LUIS EDUARDO BERTEL OLIVERA",      "1701", 48000,  1200000),
    @("CC", "1052996285", "FELIPE DE JESUS LOZANO TAFUR",     "1701", 27578,  689455),
    @("CC", "1053003136", "JUAN CARLOS MARTINEZ CAREY",       "1701", 27578,  689455),
    @("CC", "19789269",   "SAMUEL ANTONIO VANEGA SIERRA",     "1701", 27578,  689455),
    @("CC", "1083458755", "WILDIS JOSE MELENDREZ MANGA",      "1701", 27578,  737717),
    @("CC", "1053003113", "DIEGO JOSE VANEGAS CAREY",         "1701", 27578,  689455),
    @("CC", "1047371172", "KELLY DEL CARMEN BERTEL OLIVERA",  "1702", 29480,  737000),
    @("CC", "1052971182", "LUIS EDUARDO BERTEL OLIVERA",      "1702", 48000,  1200000),
    @("CC", "1047371172", "KELLY DEL CARMEN BERTEL OLIVERA",  "1703", 29480,  737000),
    @("CC", "1052971182", "LUIS EDUARDO BERTEL OLIVERA",      "1703", 48000,  1200000)
)

$row = 16
foreach ($rec in $data) {
    $ws.Cells.Item($row, 2).Value = $rec[0]

    $ws.Cells.Item($row, 3).NumberFormat = "@"
    $ws.Cells.Item($row, 3).Value = $rec[1]

    $ws.Cells.Item($row, 4).Value = $rec[2]

    $ws.Cells.Item($row, 5).NumberFormat = "@"
    $ws.Cells.Item($row, 5).Value = $rec[3]

    $ws.Cells.Item($row, 6).Value = $rec[4]
    $ws.Cells.Item($row, 7).Value = $rec[5]
    $row = $row + 1
}

# ---------------------------------------------------------------------------
# 4. Remove the now-obsolete trailing data rows (the table shrank from 23 to
#    14 rows). Deleting rows 30:38 shifts the blank spacer rows and the
#    signature footer rows up automatically, keeping merged cells correct.
# ---------------------------------------------------------------------------
$ws.Rows("30:38").Delete()

# ---------------------------------------------------------------------------
# 5. Re-fit the data columns now that the content has changed.
# ---------------------------------------------------------------------------
$ws.Columns("B:J").AutoFit()
